$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume table with refreshed figures.
# Cells whose new value looks like a plain number (e.g. "0.615") are
# written with a leading apostrophe so Excel keeps storing them as text
# (matching the source workbook, where every D/E/B/C cell is inline text),
# then the style is reset to "Normal" so no stray quote-prefix formatting
# is left behind on the cell.

$ws.Range('D2').Value = '60.719.50'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '2.647.52'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.Value = "'" + '572.90'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.69%  '
$c = $ws.Range('D6')
$c.Value = "'" + '145.58'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +2.54%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +0.44%  '
$ws.Range('E11').Value = '  +3.06%  '
$ws.Range('E12').Value = '  +2.90%  '
$ws.Range('D13').Value = '3.117.57'
$ws.Range('E13').Value = '  +1.77%  '
$c = $ws.Range('D14')
$c.Value = "'" + '25.81'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +10.90%  '
$ws.Range('D15').Value = '60.701.13'
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('E16').Value = '  +2.06%  '
$ws.Range('D17').Value = '2.658.63'
$ws.Range('E17').Value = '  +1.45%  '
$c = $ws.Range('D18')
$c.Value = "'" + '11.55'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.05%  '
$c = $ws.Range('D19')
$c.Value = "'" + '4.73'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +1.67%  '
$c = $ws.Range('D20')
$c.Value = "'" + '350.69'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +1.35%  '
$c = $ws.Range('D21')
$c.Value = "'" + '6.93'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('E23').Value = '  +0.48%  '
$c = $ws.Range('D24')
$c.Value = "'" + '63.90'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('E26').Value = '  +2.55%  '
$c = $ws.Range('D27')
$c.Value = "'" + '8.15'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +6.05%  '
$ws.Range('E28').Value = '  +10.40%  '
$ws.Range('D29').Value = '0.0₃0810'
$ws.Range('E29').Value = '  +3.48%  '
$c = $ws.Range('D30')
$c.Value = "'" + '6.67'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +4.48%  '
$c = $ws.Range('D31')
$c.Value = "'" + '169.48'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +5.69%  '
$c = $ws.Range('D33')
$c.Value = "'" + '19.72'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.26%  '
$c = $ws.Range('D34')
$c.Value = "'" + '1.08'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +11.45%  '
$c = $ws.Range('D35')
$c.Value = "'" + '4.46'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +5.63%  '
$ws.Range('E36').Value = '  +8.09%  '
$ws.Range('E37').Value = '  +2.64%  '
$c = $ws.Range('D38')
$c.Value = "'" + '330.52'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +12.35%  '
$ws.Range('E39').Value = '  +5.02%  '
$c = $ws.Range('D40')
$c.Value = "'" + '38.32'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +1.56%  '
$c = $ws.Range('D41')
$c.Value = "'" + '0.878'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +3.61%  '
$c = $ws.Range('D42')
$c.Value = "'" + '5.20'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +6.28%  '
$c = $ws.Range('D43')
$c.Value = "'" + '20.96'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +5.80%  '
$ws.Range('E44').Value = '  +3.27%  '
$c = $ws.Range('D45')
$c.Value = "'" + '133.77'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -3.63%  '
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D47')
$c.Value = "'" + '0.0560'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +2.91%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D48')
$c.Value = "'" + '0.615'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +1.55%  '
$c = $ws.Range('D49')
$c.Value = "'" + '0.998'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.18%  '
$c = $ws.Range('D50')
$c.Value = "'" + '0.0246'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +2.75%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.089.42'
$ws.Range('E51').Value = '  +3.20%  '
